# Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fallo, -1
$ws.Range("G2").Value = "Fallo"
$ws.Range("H2").Value = -1

# Row 3: Acierto, 0.62
$ws.Range("G3").Value = "Acierto"
$ws.Range("H3").Value = 0.62

# Row 6: Acierto, 1.25
$ws.Range("G6").Value = "Acierto"
$ws.Range("H6").Value = 1.25

# Rows 8 and 9: event_id stored as text -> convert to numeric value
$ws.Range("A8").Value = 14707155
$ws.Range("A9").Value = 14706848
